$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "la identificacion no le pernece al cliente"
$ws.Range("E2").Value = "El estado de cuenta no se encuentra en la bd."
$ws.Range("I2").Value = "12/02/2020 09:31 a. m.;12/02/2020 09:31 a. m.;"
$ws.Range("J2").Value = "00185537"
